# ASSIGNMENT6_ML.docx edit script
# 1) Merge the title runs ("ASSIGNMENT " / "6 :" / "- Naive Bayes ") into one run,
#    dropping the proofErr gramStart/gramEnd markers around "6 :".
$d = $word.ActiveDocument

$ok1 = $d.Content.Find.Execute(
    "ASSIGNMENT 6 :- Naive Bayes ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ASSIGNMENT 6 :- Naive Bayes ", 2)
Write-Output "title merge: $ok1"

# 2) Split "AIM: Assignment on Naive Bayes" into "AIM: Assignment of Naive Bayes.  ",
#    then break the paragraph into three: the AIM paragraph, a new OBJECTIVE
#    paragraph, and a paragraph that starts with three spaces before the
#    (unchanged) PREREQUISITE run.
$lbreak = [char]11
$umlautI = [char]239

$old2 = "AIM: Assignment on Naive Bayes" + $lbreak + "PREREQUISITE"
$new2 = "AIM: Assignment of Naive Bayes.  ^pOBJECTIVE: To apply the Na" + $umlautI + "ve Bayes algorithm on the Salary Dataset to classify individuals based on their demographic and professional attributes and predict their salary category.^p   ^lPREREQUISITE"

$ok2 = $d.Content.Find.Execute(
    $old2,
    $true, $false, $false, $false, $false, $true, 1, $false,
    $new2, 2)
Write-Output "aim/objective split: $ok2"
